# Add a new supplier/part row (GOST check via site) right after the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data down by inserting a new row at row 2.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row with the new item.
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "* Гайка шестигранная M20х1,5.5 покрытие цинковое, хроматированное ГОСТ 15521"

# Update the last active selection as recorded by Excel after the edit.
$ws.Range("B15").Select()
